$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Client Id (A2)
$ws.Range("A2").Value = "cqBFJ752"

# Update Candidate ID (B2) - numeric value
$ws.Range("B2").Value = 23092132

# Update User Name (C2)
$ws.Range("C2").Value = "amqfbph25"

# Update Exam Password (D2)
$ws.Range("D2").Value = "b32C%Kx#"

# Update First Name (F2)
$ws.Range("F2").Value = "JiMTvgut"

# Update Last Name (G2)
$ws.Range("G2").Value = "YDAH"
